$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 202.18182
$ws.Cells.Item(9, 9).Value = 165.5
$ws.Cells.Item(9, 11).Value = 165.5
$ws.Cells.Item(9, 13).Value = 3.5
$ws.Cells.Item(17, 8).Value = 1371.3636
$ws.Cells.Item(17, 10).Value = 1324.238
$ws.Cells.Item(17, 12).Value = 3972.714
$ws.Cells.Item(17, 14).Value = -4308.714
$ws.Cells.Item(32, 8).Value = 14289671
$ws.Cells.Item(32, 9).Value = 1400
$ws.Cells.Item(32, 10).Value = 16671050
$ws.Cells.Item(32, 11).Value = 1400
$ws.Cells.Item(32, 12).Value = 16671050
$ws.Cells.Item(32, 13).Value = -1074
$ws.Cells.Item(32, 14).Value = -16671702
$ws.Cells.Item(38, 8).Value = 716.2857
$ws.Cells.Item(38, 9).Value = 716.2857
$ws.Cells.Item(38, 11).Value = 2148.8571
$ws.Cells.Item(38, 13).Value = -1776.8571
$ws.Cells.Item(42, 8).Value = 359.33334
$ws.Cells.Item(42, 9).Value = 359.33334
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 1078.00002
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -848.0000199999999
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3131.647
$ws.Cells.Item(86, 10).Value = 3793.4285
$ws.Cells.Item(86, 12).Value = 3793.4285
$ws.Cells.Item(86, 14).Value = -6039.4285
$ws.Cells.Item(89, 8).Value = 3131.647
$ws.Cells.Item(89, 10).Value = 3793.4285
$ws.Cells.Item(89, 12).Value = 18967.1425
$ws.Cells.Item(89, 14).Value = -30199.1425
$ws.Cells.Item(116, 8).Value = 6749.75
$ws.Cells.Item(116, 9).Value = 5999.6665
$ws.Cells.Item(116, 11).Value = 5999.6665
$ws.Cells.Item(116, 13).Value = -2557.6665
$ws.Cells.Item(137, 8).Value = 2122.4285
$ws.Cells.Item(137, 9).Value = 1909.5135
$ws.Cells.Item(137, 10).Value = 3698
$ws.Cells.Item(137, 11).Value = 5728.5405
$ws.Cells.Item(137, 12).Value = 11094
$ws.Cells.Item(137, 13).Value = -3178.5405
$ws.Cells.Item(137, 14).Value = -16194
$ws.Cells.Item(141, 8).Value = 6210.8335
$ws.Cells.Item(141, 9).Value = 2812.3809
$ws.Cells.Item(141, 11).Value = 8437.1427
$ws.Cells.Item(141, 13).Value = -3257.1427

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1634.9062
$ws.Cells.Item(74, 9).Value = 1734.1034
$ws.Cells.Item(74, 10).Value = 676
$ws.Cells.Item(74, 11).Value = 1734.1034
$ws.Cells.Item(74, 12).Value = 676
$ws.Cells.Item(74, 13).Value = -860.1034
$ws.Cells.Item(74, 14).Value = -2424
$ws.Cells.Item(77, 8).Value = 1634.9062
$ws.Cells.Item(77, 9).Value = 1734.1034
$ws.Cells.Item(77, 10).Value = 676
$ws.Cells.Item(77, 11).Value = 8670.517
$ws.Cells.Item(77, 12).Value = 3380
$ws.Cells.Item(77, 13).Value = -4302.517
$ws.Cells.Item(77, 14).Value = -12116
$ws.Cells.Item(97, 8).Value = 1005.9
$ws.Cells.Item(97, 9).Value = 1037.1428
$ws.Cells.Item(97, 10).Value = 933
$ws.Cells.Item(97, 11).Value = 1037.1428
$ws.Cells.Item(97, 12).Value = 933
$ws.Cells.Item(97, 13).Value = -541.1428000000001
$ws.Cells.Item(97, 14).Value = -1925
$ws.Cells.Item(132, 8).Value = 1531.3112
$ws.Cells.Item(132, 10).Value = 999.3333
$ws.Cells.Item(132, 12).Value = 2997.9999
$ws.Cells.Item(132, 14).Value = -8057.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2374.2903
$ws.Cells.Item(86, 10).Value = 2857.5386
$ws.Cells.Item(86, 12).Value = 2857.5386
$ws.Cells.Item(86, 14).Value = -5103.5386
$ws.Cells.Item(89, 8).Value = 2374.2903
$ws.Cells.Item(89, 10).Value = 2857.5386
$ws.Cells.Item(89, 12).Value = 14287.693
$ws.Cells.Item(89, 14).Value = -25519.693
$ws.Cells.Item(94, 8).Value = 2421.8965
$ws.Cells.Item(94, 9).Value = 2810.4285
$ws.Cells.Item(94, 10).Value = 1402
$ws.Cells.Item(94, 11).Value = 2810.4285
$ws.Cells.Item(94, 12).Value = 1402
$ws.Cells.Item(94, 13).Value = -2359.4285
$ws.Cells.Item(94, 14).Value = -2304
$ws.Cells.Item(105, 8).Value = 2687.8
$ws.Cells.Item(105, 9).Value = 1484.3636
$ws.Cells.Item(105, 10).Value = 5997.25
$ws.Cells.Item(105, 11).Value = 1484.3636
$ws.Cells.Item(105, 12).Value = 5997.25
$ws.Cells.Item(105, 13).Value = 262.6364000000001
$ws.Cells.Item(105, 14).Value = -9491.25
$ws.Cells.Item(134, 8).Value = 2054.9688
$ws.Cells.Item(134, 9).Value = 1715.862
$ws.Cells.Item(134, 10).Value = 5333
$ws.Cells.Item(134, 11).Value = 5147.586
$ws.Cells.Item(134, 12).Value = 15999
$ws.Cells.Item(134, 13).Value = -2612.586
$ws.Cells.Item(134, 14).Value = -21069

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1884.6
$ws.Cells.Item(31, 9).Value = 1793.8572
$ws.Cells.Item(31, 11).Value = 1793.8572
$ws.Cells.Item(31, 13).Value = -1498.8572
$ws.Cells.Item(34, 8).Value = 1884.6
$ws.Cells.Item(34, 9).Value = 1793.8572
$ws.Cells.Item(34, 11).Value = 1793.8572
$ws.Cells.Item(34, 13).Value = -1591.8572
$ws.Cells.Item(107, 8).Value = 14572.333
$ws.Cells.Item(107, 9).Value = 1431.6666
$ws.Cells.Item(107, 10).Value = 17857.5
$ws.Cells.Item(107, 11).Value = 1431.6666
$ws.Cells.Item(107, 12).Value = 17857.5
$ws.Cells.Item(107, 13).Value = 488.3334
$ws.Cells.Item(107, 14).Value = -21697.5
$ws.Cells.Item(132, 8).Value = 3166.7778
$ws.Cells.Item(132, 9).Value = 2714.4285
$ws.Cells.Item(132, 11).Value = 8143.2855
$ws.Cells.Item(132, 13).Value = -5613.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 695.125
$ws.Cells.Item(26, 9).Value = 73.583336
$ws.Cells.Item(26, 10).Value = 2559.75
$ws.Cells.Item(26, 11).Value = 220.750008
$ws.Cells.Item(26, 12).Value = 7679.25
$ws.Cells.Item(26, 13).Value = 67.24999199999999
$ws.Cells.Item(26, 14).Value = -8255.25
$ws.Cells.Item(69, 8).Value = 3999.75
$ws.Cells.Item(69, 9).Value = 2500
$ws.Cells.Item(69, 11).Value = 7500
$ws.Cells.Item(69, 13).Value = -6689
$ws.Cells.Item(72, 8).Value = 3999.75
$ws.Cells.Item(72, 9).Value = 2500
$ws.Cells.Item(72, 11).Value = 22500
$ws.Cells.Item(72, 13).Value = -18444
$ws.Cells.Item(113, 8).Value = 1437.7858
$ws.Cells.Item(113, 10).Value = 1552.5714
$ws.Cells.Item(113, 12).Value = 4657.7142
$ws.Cells.Item(113, 14).Value = -8997.7142
$ws.Cells.Item(138, 8).Value = 5592.76
$ws.Cells.Item(138, 9).Value = 3267.182
$ws.Cells.Item(138, 11).Value = 9801.545999999998
$ws.Cells.Item(138, 13).Value = -4661.545999999998
$ws.Cells.Item(140, 8).Value = 10873531
$ws.Cells.Item(140, 10).Value = 7333.3335
$ws.Cells.Item(140, 12).Value = 22000.0005
$ws.Cells.Item(140, 14).Value = -32360.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 6840.8
$ws.Cells.Item(126, 9).Value = 7700
$ws.Cells.Item(126, 10).Value = 6268
$ws.Cells.Item(126, 11).Value = 23100
$ws.Cells.Item(126, 12).Value = 18804
$ws.Cells.Item(126, 13).Value = -20630
$ws.Cells.Item(126, 14).Value = -23744
$ws.Cells.Item(132, 8).Value = 2426.5
$ws.Cells.Item(132, 9).Value = 2287.4285
$ws.Cells.Item(132, 10).Value = 3400
$ws.Cells.Item(132, 11).Value = 6862.2855
$ws.Cells.Item(132, 12).Value = 10200
$ws.Cells.Item(132, 13).Value = -4332.2855
$ws.Cells.Item(132, 14).Value = -15260

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6432.696
$ws.Cells.Item(40, 10).Value = 3688.6
$ws.Cells.Item(40, 12).Value = 3688.6
$ws.Cells.Item(40, 14).Value = -3960.6
$ws.Cells.Item(46, 8).Value = 2991
$ws.Cells.Item(46, 10).Value = 2991
$ws.Cells.Item(46, 12).Value = 2991
$ws.Cells.Item(46, 14).Value = -3367
$ws.Cells.Item(136, 8).Value = 2938.9333
$ws.Cells.Item(136, 9).Value = 2969.1853
$ws.Cells.Item(136, 11).Value = 8907.555899999999
$ws.Cells.Item(136, 13).Value = -6357.555899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 10909.667
$ws.Cells.Item(45, 9).Value = 9284
$ws.Cells.Item(45, 10).Value = 11722.5
$ws.Cells.Item(45, 11).Value = 9284
$ws.Cells.Item(45, 12).Value = 11722.5
$ws.Cells.Item(45, 13).Value = -8793
$ws.Cells.Item(45, 14).Value = -12704.5
$ws.Cells.Item(96, 8).Value = 37268
$ws.Cells.Item(96, 10).Value = 4542
$ws.Cells.Item(96, 12).Value = 4542
$ws.Cells.Item(96, 14).Value = -7288
$ws.Cells.Item(136, 8).Value = 2014.5
$ws.Cells.Item(136, 9).Value = 1151.5264
$ws.Cells.Item(136, 11).Value = 3454.5792
$ws.Cells.Item(136, 13).Value = -904.5792000000001
